# Remove the "Keep SPRING WG in the loop for SR aspects" bullet from the
# "Next Steps" slide (slide 12), per the commit diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)              # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange

# Find the paragraph with the exact text and delete it (removes the whole
# paragraph, including its trailing line break, shifting later ones up).
for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.TrimEnd("`r", "`n") -eq "Keep SPRING WG in the loop for SR aspects") {
        $para.Delete()
    }
}
